$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "A: What is Github" paragraph: drop the gramStart/gramEnd proofErr
#    markers that used to bracket the run(s) spelling "Github".
#    iron_native auto-manages proofErr markers based on spelling/grammar
#    state, so flip the checked flags to make it drop stale markers.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(8).Range.GrammarChecked = $true

# ---------------------------------------------------------------------
# 2) Empty "Workflow" heading paragraph, just above "Commit/Directory/
#    Clone/<empty>" list and right before the "Questions:" paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(37).Range.InsertBefore("Workflow")

# ---------------------------------------------------------------------
# 3) "Have you used a project on Github?" -> "Do you know which is the
#    most used language on GitHub?"
# ---------------------------------------------------------------------
$q1 = $d.Paragraphs.Item(39).Range
$q1.Text = "Do you know which is the most used language on GitHub?"

# ---------------------------------------------------------------------
# 4) "Do you know the difference between saving a file and making a
#    commit?" -> "The difference between saving a file and making a
#    commit?" (capitalized T split into its own run)
# ---------------------------------------------------------------------
$q2 = $d.Paragraphs.Item(40).Range
$q2.Text = "The difference between saving a file and making a commit?"
